$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.634.46"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.743.88"
$ws.Range("E3").Value = "  +4.52%  "
$ws.Range("D5").Value = "'611.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.59%  "
$ws.Range("D6").Value = "'177.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("D7").Value = "3.740.10"
$ws.Range("E7").Value = "  +4.40%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("D11").Value = "'6.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "'40.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.30%  "
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "4.365.79"
$ws.Range("E15").Value = "  +5.35%  "
$ws.Range("D16").Value = "3.746.56"
$ws.Range("E16").Value = "  +4.55%  "
$ws.Range("D17").Value = "69.661.59"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").Value = "'512.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "'9.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").Value = "'0.722"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.51%  "
$ws.Range("D24").Value = "'87.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  +5.02%  "
$ws.Range("D26").Value = "'13.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").Value = "'11.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "'0.0000128"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.17%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "'2.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("D32").Value = "'7.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").Value = "'31.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "'0.114"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "'6.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").Value = "'0.336"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").Value = "'0.131"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.15%  "
$ws.Range("D41").Value = "'51.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'44.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.063.91"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'416.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.76%  "
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'135.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
